$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new column before DK, containing 08-nov header + "-" values ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DK1").EntireColumn.Insert()

$ws1.Range("DK1").Value = "08-nov"
$ws1.Range("DK2:DK25").Value = "-"

# --- Sheet "Gaz": append row 144 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A144").Value = "'2025-11-06"
$wsGaz.Range("A144").Style = "Normal"
$wsGaz.Range("B144").Value = 30.35

# --- Sheet "CO2": append row 144 ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A144").Value = "'2025-11-06"
$wsCO2.Range("A144").Style = "Normal"
$wsCO2.Range("B144").Value = 79.94
